$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-06-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-07 Saturday", 2) | Out-Null

# Replacement values for each cell of the 20x5 answers table,
# addressed by (row, column) position rather than by text search,
# since several old answers are not unique within the table.
$answers = @(
    @("7+89=96", "18+26=44", "37+38=75", "90-71=19", "61-35=26"),
    @("38+46=84", "49+42=91", "53-28=25", "78+9=87", "36-7=29"),
    @("18+44=62", "64-8=56", "56+39=95", "61-35=26", "45+7=52"),
    @("6+55=61", "78+5=83", "74+8=82", "15+66=81", "50-45=5"),
    @("8+14=22", "60-24=36", "8+78=86", "42+49=91", "80-62=18"),
    @("7+35=42", "31-5=26", "17+35=52", "19+16=35", "70-55=15"),
    @("16+37=53", "61-12=49", "41-24=17", "3+29=32", "30-5=25"),
    @("92-39=53", "34-26=8", "67+27=94", "81-8=73", "92-63=29"),
    @("58+38=96", "81-68=13", "71-8=63", "27+54=81", "68+27=95"),
    @("16+26=42", "8+17=25", "91-2=89", "82-29=53", "48+9=57"),
    @("19+38=57", "6+78=84", "64-56=8", "9+19=28", "41-12=29"),
    @("19+79=98", "80-9=71", "84-75=9", "36+56=92", "26+39=65"),
    @("44+17=61", "90-82=8", "88-29=59", "90-9=81", "60-3=57"),
    @("17+4=21", "19+25=44", "60-48=12", "93-6=87", "74-37=37"),
    @("76-48=28", "45-17=28", "96-87=9", "8+36=44", "28+54=82"),
    @("78+13=91", "49+19=68", "69+29=98", "90-81=9", "17+7=24"),
    @("16+65=81", "83-4=79", "92-27=65", "84-58=26", "15+57=72"),
    @("19+34=53", "5+87=92", "19+65=84", "52-23=29", "6+25=31"),
    @("73-55=18", "48+19=67", "19+12=31", "20-17=3", "15+37=52"),
    @("50-41=9", "52-49=3", "28+63=91", "88-19=69", "48+25=73")
)

$table = $d.Tables.Item(1)
for ($r = 0; $r -lt $answers.Length; $r++) {
    $row = $answers[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $table.Cell($r + 1, $c + 1)
        $cell.Range.Text = $row[$c]
    }
}

Write-Output "done"
